$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 282. Excel shifts rows 282:372 down to
# 283:373 (content + formatting), which reproduces the "row N gets the
# data that used to live in row N-1" shift seen across the whole block,
# and pushes the former last row (372) down into the newly-created 373.
$ws.Rows("282:282").Insert()

# Populate the newly blank row 282 with the new weekly data point.
$ws.Range("A282").Value = 3
$ws.Range("B282").Value = "Femacal de La Calera"
$ws.Range("C282").Value = "Coquimbo"
$ws.Range("D282").Value = 44463
$ws.Range("E282").Value = 5
$ws.Range("F282").Value = 100112024
$ws.Range("G282").Value = "Choclo"
$ws.Range("H282").Value = "Dulce o Americano"
$ws.Range("I282").Value = "Primera"
$ws.Range("J282").Value = 40
$ws.Range("K282").Value = 32000
$ws.Range("L282").Value = 32000
$ws.Range("M282").Value = 32000
$ws.Range("N282").Value = "`$/malla 70 unidades"
$ws.Range("O282").Value = "Región de Arica y Parinacota"
$ws.Range("P282").Value = 457
$ws.Range("Q282").Value = 70
$ws.Range("R282").Value = "Hortaliza"
